$wb = $excel.ActiveWorkbook

# Add the new "Color Scheme Features" sheet at the very end of the workbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Color Scheme Features"

# Write the new unique strings in the same order they were authored so the
# shared string table indices line up with the target workbook.
$ws.Range("C1").Value = "Feature Title"
$ws.Range("D1").Value = "Feature Price"
$ws.Range("G1").Value = "Material"
$ws.Range("H1").Value = "Manufacturer"
$ws.Range("I1").Value = "Name"
$ws.Range("K1").Value = "Feature Image"
$ws.Range("J1").Value = "Manufacturer Id"
$ws.Range("E1").Value = "(Upgrade=1, Base=0)"
$ws.Range("C2").Value = "Windows"
$ws.Range("H2").Value = "ML"
$ws.Range("F1").Value = "Upgraded Type(Concrete=1, Window=2, Wall=3), Base=0"
$ws.Range("G2").Value = "windows"
$ws.Range("I2").Value = "Black"
$ws.Range("J2").Value = "N/A"
$ws.Range("K2").Value = "elante-elante-basic-windows.png"
$ws.Range("A1").Value = "Elevation or Elevation Type Title"

# Cells that reuse already-existing shared strings
$ws.Range("B1").Value = "Color Scheme Title"
$ws.Range("A2").Value = "Elante"
$ws.Range("B2").Value = "Elante Basic"

# Numeric data cells
$ws.Range("D2").Value = 500
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 2

# Bold header row (matches style index 2 used on all other sheets' headers)
$ws.Range("A1:K1").Font.Bold = $true

# Column widths (bestFit custom widths matching the sample workbook)
$ws.Columns.Item(1).ColumnWidth = 28.21875
$ws.Columns.Item(2).ColumnWidth = 16.6640625
$ws.Columns.Item(3).ColumnWidth = 11.33203125
$ws.Columns.Item(4).ColumnWidth = 11.88671875
$ws.Columns.Item(5).ColumnWidth = 18.21875
$ws.Columns.Item(6).ColumnWidth = 49.33203125
$ws.Columns.Item(7).ColumnWidth = 8.21875
$ws.Columns.Item(8).ColumnWidth = 12.6640625
$ws.Columns.Item(9).ColumnWidth = 6
$ws.Columns.Item(10).ColumnWidth = 14.88671875
$ws.Columns.Item(11).ColumnWidth = 28.44140625

$ws.Range("B2").Select()

# Workbook-level view tweaks
$excel.ActiveWindow.DisplayedSheets = [System.Reflection.Missing]::Value
$wb.Windows.Item(1).ScrollWorkbookTabs(2)

# Selection tweak on "Elevation Color Schemes"
$wb.Worksheets.Item("Elevation Color Schemes").Range("D3").Select()

# Re-select the new sheet so it is the active tab again
$ws.Activate()
$ws.Range("B2").Select()
